$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2408311724041815
$ws.Range("C2").Value = 0.2408311724041815
$ws.Range("D2").Value = 0.2408311724041816

$ws.Range("B3").Value = 0.3742737492009104
$ws.Range("C3").Value = 0.3742737492009104
$ws.Range("D3").Value = 0.3742737492009104

$ws.Range("B4").Value = 0.384895078394908
$ws.Range("C4").Value = 0.3848950783949081
$ws.Range("D4").Value = 0.3848950783949081

$ws.Range("E5").Value = 0.2231518276835295
$ws.Range("F5").Value = 0.2231518276835295
$ws.Range("G5").Value = 0.2231518276835295

$ws.Range("E6").Value = 0.3719475108519568
$ws.Range("F6").Value = 0.3719475108519568
$ws.Range("G6").Value = 0.3719475108519568

$ws.Range("E7").Value = 0.4049006614645136
$ws.Range("F7").Value = 0.4049006614645137
$ws.Range("G7").Value = 0.4049006614645136
